# Commiting code that gets test data from excel files
#
# The workbook originally contains a single "Language" sheet.
# This script adds three new sheets (Skill, Education, Certification)
# before the existing "Language" sheet, populates them with test data,
# sizes their columns, and restores per-sheet selections - matching the
# layout used by the SpecFlow test data workbook.

$wb = $excel.ActiveWorkbook

# Give the workbook its VBA project code name (cosmetic / best effort).
$wb.CodeName = "ThisWorkbook"

# --- Create the new sheets, in order, right after one another so the
#     final left-to-right tab order is Skill, Education, Certification,
#     Language (Language stays where it always was). -------------------
$skill = $wb.Worksheets.Add()
$skill.Name = "Skill"
$skill.CodeName = "Sheet1"

$education = $wb.Worksheets.Add($null, $skill)
$education.Name = "Education"
$education.CodeName = "Sheet2"

$certification = $wb.Worksheets.Add($null, $education)
$certification.Name = "Certification"
$certification.CodeName = "Sheet3"

# Look the original sheet back up by name *after* the new sheets have
# been inserted ahead of it, so this reference resolves to "Language"
# rather than whatever now sits at its old position.
$lang = $wb.Worksheets.Item("Language")
$lang.CodeName = "Sheet4"

# --- Skill sheet ------------------------------------------------------
$skill.Range("A1").Value = "Skill"
$skill.Range("B1").Value = "Level"
$skill.Range("A2").Value = "Skill1"
$skill.Range("B2").Value = "Intermediate"

$skill.Columns("B").ColumnWidth = 11.88

# --- Education sheet ----------------------------------------------------
# Values are written in this particular order so that freshly introduced
# strings land in the workbook's shared-string table in the same order
# they originally were typed in (Country, Title, Degree, AUT,
# New Zealand, Institute, Associate, Test Analyst, Year).
$education.Range("A1").Value = "Country"
$education.Range("C1").Value = "Title"
$education.Range("D1").Value = "Degree"
$education.Range("B2").Value = "AUT"
$education.Range("A2").Value = "New Zealand"
$education.Range("B1").Value = "Institute"
$education.Range("C2").Value = "Associate"
$education.Range("D2").Value = "Test Analyst"
$education.Range("E1").Value = "Year"
$education.Range("E2").Value = 2022

$education.Columns("A").ColumnWidth = 11.74
$education.Columns("B").ColumnWidth = 7.74
$education.Columns("C").ColumnWidth = 8.59
$education.Columns("D").ColumnWidth = 10.88
$education.Columns("E").ColumnWidth = 4.17

# --- Certification sheet ----------------------------------------------
$certification.Range("A1").Value = "Certificate"
$certification.Range("B1").Value = "From"
$certification.Range("C1").Value = "Year"
$certification.Range("A2").Value = "Certified Tester Foundation Level"
$certification.Range("B2").Value = "ISTQB"
$certification.Range("C2").Value = 2022

$certification.Columns("A").ColumnWidth = 30.45

# --- Selections ---------------------------------------------------------
# Select ranges in this order so the last-selected sheet (Certification)
# ends up being the active / visible tab, matching activeTab="2".
$skill.Range("F13").Select()
$education.Range("I12").Select()
$lang.Range("F12").Select()
$certification.Range("C3").Select()
